$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 253, shifting existing rows 253:320 down to 254:321.
$ws.Rows.Item(253).Insert()

# Populate the newly inserted row 253 with the new "Navel Late" price entry.
$ws.Range("A253").Value = 11
$ws.Range("B253").Value = "Vega Monumental Concepción"
$ws.Range("C253").Value = "Bíobío"
$ws.Range("D253").Value = 44782
$ws.Range("E253").Value = 8
$ws.Range("F253").Value = "Fruta"
$ws.Range("G253").Value = 100102
$ws.Range("H253").Value = "Cítricos"
$ws.Range("I253").Value = 100102005
$ws.Range("J253").Value = "Naranja"
$ws.Range("K253").Value = "Navel Late"
$ws.Range("L253").Value = "Primera"
$ws.Range("M253").Value = 350
$ws.Range("N253").Value = 6000
$ws.Range("O253").Value = 6500
$ws.Range("P253").Value = 6214
$ws.Range("Q253").Value = "$/bandeja 15 kilos granel"
$ws.Range("R253").Value = "Región de O'Higgins"
$ws.Range("S253").Value = 414
$ws.Range("T253").Value = 15

# Match the date-number format already used by the other rows' Fecha column.
$ws.Range("D253").NumberFormat = $ws.Range("D254").NumberFormat
